# adding utility to figure
#
# Applies the edits described by the diff:
#  - add a new "Super confusing, value versus impact…" note in A3
#  - rename the "rating_value" header to "rating_text"
#  - reverse + relabel the five rating-bucket rows (A7:A11, A12:A16,
#    A17:A21, A22:A26) from "<level> value" to bare "<level>", in
#    reverse order, and clear the special border/shading formatting
#    that used to highlight the first/last rows of each block
#  - move the active selection to B32

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New explanatory note in the previously-empty row 3
$ws.Range("A3").Value = "Super confusing, value versus impact…"

# Header rename
$ws.Range("A6").Value = "rating_text"

# Each block of 5 rows (one per confidence column) used to read, top to
# bottom: "very high value", "high value", "medium value", "low value",
# "very low value". They now read, top to bottom: "very low", "low",
# "medium", "high", "very high" -- and the per-row emphasis formatting
# (borders/bold first & last row of each block) is cleared.
$blockStarts = @(7, 12, 17, 22)
$labels = @("very low", "low", "medium", "high", "very high")

foreach ($start in $blockStarts) {
    for ($i = 0; $i -lt 5; $i++) {
        $row = $start + $i
        $cell = $ws.Range("A" + $row)
        $cell.Value = $labels[$i]
        $cell.ClearFormats()
    }
}

# Selection moved to B32
$ws.Range("B32").Select() | Out-Null
